$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.370.31"
$ws.Range("E2").Value = "  -3.77%  "

$ws.Range("D3").Value = "3.552.42"
$ws.Range("E3").Value = "  -4.48%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.13"
$ws.Range("E5").Value = "  -7.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.72"
$ws.Range("E6").Value = "  -4.28%  "

$ws.Range("D7").Value = "3.550.19"
$ws.Range("E7").Value = "  -4.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  -4.18%  "

$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.669"
$ws.Range("E10").Value = "  -7.73%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("E11").Value = "  -8.02%  "

$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.45"
$ws.Range("E12").Value = "  -7.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -10.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.77"
$ws.Range("E14").Value = "  -6.52%  "

$ws.Range("D15").Value = "4.121.79"
$ws.Range("E15").Value = "  -4.48%  "

$ws.Range("D16").Value = "3.564.79"
$ws.Range("E16").Value = "  -4.28%  "

$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("E18").Value = "  -6.59%  "

$ws.Range("D19").Value = "66.407.31"
$ws.Range("E19").Value = "  -3.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.02"
$ws.Range("E20").Value = "  -7.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("E21").Value = "  -8.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.31"
$ws.Range("E22").Value = "  -6.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("E23").Value = "  -9.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.16"
$ws.Range("E24").Value = "  -5.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.12"
$ws.Range("E25").Value = "  -2.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.90"
$ws.Range("E26").Value = "  -6.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.34"
$ws.Range("E27").Value = "  -6.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.05"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.54"
$ws.Range("E29").Value = "  -7.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.79"
$ws.Range("E30").Value = "  -9.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.50"
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.72"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "628.04"
$ws.Range("E33").Value = "  -1.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.10"
$ws.Range("E34").Value = "  -5.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.113"
$ws.Range("E35").Value = "  -8.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.13"
$ws.Range("E36").Value = "  -6.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.63"
$ws.Range("E37").Value = "  -10.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.400"
$ws.Range("E38").Value = "  -3.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "0.0₃0743"
$ws.Range("E40").Value = "  -10.60%  "

$ws.Range("D41").Value = "3.134.81"
$ws.Range("E41").Value = "  +7.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.132"
$ws.Range("E42").Value = "  -6.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.92"
$ws.Range("E44").Value = "  -4.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.62"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0409"
$ws.Range("E46").Value = "  -9.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.129"
$ws.Range("E47").Value = "  -7.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.05"
$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.02"
$ws.Range("E49").Value = "  -4.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.38"
$ws.Range("E50").Value = "  -10.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.72"
$ws.Range("E51").Value = "  -2.23%  "
